$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update cell C2: "generarMatrix()" -> "generateMatrix()"
$ws.Range("C2").Value = "generateMatrix()"

# Update cell A2: "... y no más" -> "... y ceros de ahí para alla"
$ws.Range("A2").Value = "Comprobar que el metodo agrega los numeros hasta el número que entra por parametro y ceros de ahí para alla"

# Update selection to A3
$ws.Range("A3").Select()
